$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "a0"
$ws.Range("A3").Value = "a1"
$ws.Range("A4").Value = "a2"
$ws.Range("A5").Value = "b0"
$ws.Range("A6").Value = "b1"
$ws.Range("A7").Value = "c0"
$ws.Range("A8").Value = "c1"
$ws.Range("A9").Value = "c2"
$ws.Range("A10").Value = "c3"
